$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr1 = New-Object 'object[,]' 24,7
$arr1[0,0] = 2.945272810113571
$arr1[0,1] = 2.857492738202937
$arr1[0,2] = 9.660890158981333
$arr1[0,3] = 19.1194699430181
$arr1[0,4] = 21.81094436702325
$arr1[0,5] = 11.33996407602222
$arr1[0,6] = 15.88392123655433
$arr1[1,0] = 2.892394151407728
$arr1[1,1] = 2.844948643225307
$arr1[1,2] = 9.848614818674685
$arr1[1,3] = 18.79176649028076
$arr1[1,4] = 21.06375382383404
$arr1[1,5] = 11.31536942291477
$arr1[1,6] = 15.71536822543159
$arr1[2,0] = 2.860097505996259
$arr1[2,1] = 2.837521751649239
$arr1[2,2] = 9.969416308665235
$arr1[2,3] = 18.59435496015756
$arr1[2,4] = 20.59974027787345
$arr1[2,5] = 11.30304127565922
$arr1[2,6] = 15.61609346715682
$arr1[3,0] = 2.846997967653003
$arr1[3,1] = 2.834566522608437
$arr1[3,2] = 10.02003913748257
$arr1[3,3] = 18.51498364911185
$arr1[3,4] = 20.40969295833277
$arr1[3,5] = 11.29871797271778
$arr1[3,6] = 15.57674582186314
$arr1[4,0] = 2.844827066512321
$arr1[4,1] = 2.834080174129344
$arr1[4,2] = 10.02852937122893
$arr1[4,3] = 18.50187253646626
$arr1[4,4] = 20.37808886344588
$arr1[4,5] = 11.29804248106588
$arr1[4,6] = 15.57028039933719
$arr1[5,0] = 2.859920566905787
$arr1[5,1] = 2.837481605044526
$arr1[5,2] = 9.970093372956528
$arr1[5,3] = 18.59328001889471
$arr1[5,4] = 20.59718061962194
$arr1[5,5] = 11.30298013003119
$arr1[5,6] = 15.61555826601833
$arr1[6,0] = 2.92701583801165
$arr1[6,1] = 2.853111678337652
$arr1[6,2] = 9.724469861591388
$arr1[6,3] = 19.0057615609504
$arr1[6,4] = 21.55460664029507
$arr1[6,5] = 11.33091041941717
$arr1[6,6] = 15.82495784515375
$arr1[7,0] = 3.059142478383146
$arr1[7,1] = 2.885851857482808
$arr1[7,2] = 9.286621938734921
$arr1[7,3] = 19.83947044754661
$arr1[7,4] = 23.37618954111514
$arr1[7,5] = 11.40751328518121
$arr1[7,6] = 16.26686793387874
$arr1[8,0] = 3.155504498033453
$arr1[8,1] = 2.911049996021791
$arr1[8,2] = 8.991485936072337
$arr1[8,3] = 20.46005442701366
$arr1[8,4] = 24.66312781706535
$arr1[8,5] = 11.47682912063841
$arr1[8,6] = 16.60754844449703
$arr1[9,0] = 3.198976894209952
$arr1[9,1] = 2.92273283783765
$arr1[9,2] = 8.862956075684419
$arr1[9,3] = 20.74266320019359
$arr1[9,4] = 25.23434298857242
$arr1[9,5] = 11.51112086403216
$arr1[9,6] = 16.76533313915579
$arr1[10,0] = 3.215369954508149
$arr1[10,1] = 2.927185959140635
$arr1[10,2] = 8.815107336282885
$arr1[10,3] = 20.84960621954469
$arr1[10,4] = 25.44838181166293
$arr1[10,5] = 11.52449591266434
$arr1[10,6] = 16.82542842665081
$arr1[11,0] = 3.211842742219154
$arr1[11,1] = 2.926225647562501
$arr1[11,2] = 8.825375841125439
$arr1[11,3] = 20.82657928611917
$arr1[11,4] = 25.40238887410828
$arr1[11,5] = 11.52159815628518
$arr1[11,6] = 16.81247138206119
$arr1[12,0] = 3.200327036938166
$arr1[12,1] = 2.923098630768104
$arr1[12,2] = 8.859003047450321
$arr1[12,3] = 20.75146362637511
$arr1[12,4] = 25.2519987329012
$arr1[12,5] = 11.51221347513658
$arr1[12,6] = 16.77027061965267
$arr1[13,0] = 3.193263866624943
$arr1[13,1] = 2.921186953359234
$arr1[13,2] = 8.879707799169331
$arr1[13,3] = 20.70543994586988
$arr1[13,4] = 25.15957888350052
$arr1[13,5] = 11.50651558463534
$arr1[13,6] = 16.74446472541444
$arr1[14,0] = 3.152654714442388
$arr1[14,1] = 2.910290705069418
$arr1[14,2] = 9.000000845766174
$arr1[14,3] = 20.44158114784532
$arr1[14,4] = 24.62549176132642
$arr1[14,5] = 11.47464299384559
$arr1[14,6] = 16.59728837373711
$arr1[15,0] = 3.127637106863963
$arr1[15,1] = 2.903660765488691
$arr1[15,2] = 9.075263185443863
$arr1[15,3] = 20.27970415695155
$arr1[15,4] = 24.29403762348669
$arr1[15,5] = 11.45579199030456
$arr1[15,6] = 16.50767821124227
$arr1[16,0] = 3.113214410178838
$arr1[16,1] = 2.899868270947553
$arr1[16,2] = 9.119091377741729
$arr1[16,3] = 20.18663203040163
$arr1[16,4] = 24.10206721953766
$arr1[16,5] = 11.44520961949479
$arr1[16,6] = 16.45640455345792
$arr1[17,0] = 3.108325949961692
$arr1[17,1] = 2.8985878596412
$arr1[17,2] = 9.13402351202612
$arr1[17,3] = 20.15512913480406
$arr1[17,4] = 24.03684869258807
$arr1[17,5] = 11.4416715176607
$arr1[17,6] = 16.43909192640924
$arr1[18,0] = 3.130303837245287
$arr1[18,1] = 2.904364392237945
$arr1[18,2] = 9.06719557576648
$arr1[18,3] = 20.29693346276987
$arr1[18,4] = 24.32946049425547
$arr1[18,5] = 11.45777183020653
$arr1[18,6] = 16.51719007669677
$arr1[19,0] = 3.203711480657593
$arr1[19,1] = 2.924016342729327
$arr1[19,2] = 8.849103600223792
$arr1[19,3] = 20.77352986547625
$arr1[19,4] = 25.2962351676085
$arr1[19,5] = 11.51495947373993
$arr1[19,6] = 16.78265708812577
$arr1[20,0] = 3.251278300092378
$arr1[20,1] = 2.937028330377697
$arr1[20,2] = 8.711363120930281
$arr1[20,3] = 21.08452578551665
$arr1[20,4] = 25.91477293681093
$arr1[20,5] = 11.55460143452917
$arr1[20,6] = 16.95814302529789
$arr1[21,0] = 3.225933767423637
$arr1[21,1] = 2.930069044487508
$arr1[21,2] = 8.784439315334232
$arr1[21,3] = 20.91862383940697
$arr1[21,4] = 25.58593237543431
$arr1[21,5] = 11.53323896497671
$arr1[21,6] = 16.86432013052328
$arr1[22,0] = 3.129098330956484
$arr1[22,1] = 2.904046222771437
$arr1[22,2] = 9.070841200980647
$arr1[22,3] = 20.28914410592916
$arr1[22,4] = 24.31345019830551
$arr1[22,5] = 11.45687594868073
$arr1[22,6] = 16.51288899650855
$arr1[23,0] = 3.02344984620726
$arr1[23,1] = 2.8767842365738
$arr1[23,2] = 9.400397923122076
$arr1[23,3] = 19.6120089327625
$arr1[23,4] = 22.8912849762203
$arr1[23,5] = 11.38447844197752
$arr1[23,6] = 16.14428118618297
$ws.Range("C2:I25").Value = $arr1

$arr2 = New-Object 'object[,]' 24,3
$arr2[0,0] = 18.80208496337583
$arr2[0,1] = 17.24026566932926
$arr2[0,2] = 16.47549961453703
$arr2[1,0] = 18.04094021997504
$arr2[1,1] = 16.93383788956176
$arr2[1,2] = 16.30706827880249
$arr2[2,0] = 17.55742226830502
$arr2[2,1] = 16.74479773534667
$arr2[2,2] = 16.20812450743713
$arr2[3,0] = 17.35658199820171
$arr2[3,1] = 16.66763592129333
$arr2[3,2] = 16.16897483646332
$arr2[4,0] = 17.323011300689
$arr2[4,1] = 16.65481853048884
$arr2[4,2] = 16.1625460558897
$arr2[5,0] = 17.55472869087123
$arr2[5,1] = 16.74375748498235
$arr2[5,2] = 16.20759172268129
$arr2[6,0] = 18.54314704846287
$arr2[6,1] = 17.13485260326989
$arr2[6,2] = 16.41652422298732
$arr2[7,0] = 20.34311927345325
$arr2[7,1] = 17.8905349282614
$arr2[7,2] = 16.85960358093503
$arr2[8,0] = 21.56988343456625
$arr2[8,1] = 18.43350607643818
$arr2[8,2] = 17.20247180148453
$arr2[9,0] = 22.10531176182099
$arr2[9,1] = 18.67680075638491
$arr2[9,2] = 17.36154851418575
$arr2[10,0] = 22.30468236636218
$arr2[10,1] = 18.76831451019662
$arr2[10,2] = 17.42217613271754
$arr2[11,0] = 22.26189672985715
$arr2[11,1] = 18.74863409650544
$arr2[11,2] = 17.40910252134551
$arr2[12,0] = 22.12178260720155
$arr2[12,1] = 18.68434247013467
$arr2[12,2] = 17.36652892511394
$arr2[13,0] = 22.03551437637885
$arr2[13,1] = 18.64487927791798
$arr2[13,2] = 17.34050027685761
$arr2[14,0] = 21.53442495070978
$arr2[14,1] = 18.41752436090751
$arr2[14,2] = 17.1921333048274
$arr2[15,0] = 21.22112726946788
$arr2[15,1] = 18.27703929312078
$arr2[15,2] = 17.10186909703969
$arr2[16,0] = 21.03880285341849
$arr2[16,1] = 18.19589193012453
$arr2[16,2] = 17.05024717682241
$arr2[17,0] = 20.97671044280456
$arr2[17,1] = 18.16836036273692
$arr2[17,2] = 17.0328214155287
$arr2[18,0] = 21.2546991585039
$arr2[18,1] = 18.29203039991417
$arr2[18,2] = 17.11144769292883
$arr2[19,0] = 22.16303031519017
$arr2[19,1] = 18.70324385986559
$arr2[19,2] = 17.37902373282015
$arr2[20,0] = 22.73690274819445
$arr2[20,1] = 18.96836571639813
$arr2[20,2] = 17.55613914195533
$arr2[21,0] = 22.43246363937578
$arr2[21,1] = 18.82722379674622
$arr2[21,2] = 17.46142351301616
$arr2[22,0] = 21.23952816673902
$arr2[22,1] = 18.28525409966325
$arr2[22,2] = 17.10711635747137
$arr2[23,0] = 19.87226317766147
$arr2[23,1] = 17.68787599418601
$arr2[23,2] = 16.73647062420718
$ws.Range("M2:O25").Value = $arr2

Write-Output "done"